# Applies the cryptos.xlsx price/volume refresh (GitHub Actions data update)
# described by the commit diff: updates Price (D) and Volume(1h) (E) values for
# each coin row, plus the two rank swaps (Cosmos/OKB and Kaspa/dogwifhat).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.484.89"
$ws.Range("E2").Value = "  -1.86%  "

# Row 3
$ws.Range("D3").Value = "2.991.41"
$ws.Range("E3").Value = "  -1.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.45%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$ws.Range("D9").Value = "2.989.48"
$ws.Range("E9").Value = "  -1.10%  "

# Row 10
$ws.Range("E10").Value = "  -2.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.21%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.15%  "

# Row 15
$ws.Range("E15").Value = "  +2.16%  "

# Row 16
$ws.Range("D16").Value = "3.487.30"
$ws.Range("E16").Value = "  -1.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "

# Row 18
$ws.Range("D18").Value = "61.438.17"
$ws.Range("E18").Value = "  -1.97%  "

# Row 19
$ws.Range("D19").Value = "2.988.69"
$ws.Range("E19").Value = "  -1.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.08%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.83%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.48%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("E29").Value = "  +1.99%  "

# Row 30
$ws.Range("E30").Value = "  +0.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.48%  "

# Row 32
$ws.Range("E32").Value = "  -3.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "

# Row 34
$ws.Range("E34").Value = "  +0.34%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0815"
$ws.Range("E35").Value = "  +2.68%  "

# Row 36
$ws.Range("E36").Value = "  -1.44%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.69%  "

# Row 39
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.78%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.99%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "397.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0352"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "38.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("E46").Value = "  -5.13%  "

# Row 47
$ws.Range("D47").Value = "2.719.11"
$ws.Range("E47").Value = "  -3.00%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.78%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
